$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing row (22) into a new row 23, preserving the
# same cell styles (hyperlink style on column B, "cleared" style on column C)
$ws.Rows("22:22").Copy()
$ws.Rows("23:23").Insert()

# Set the new link + its description
$newUrl = "https://developer.mozilla.org/es/docs/Web/JavaScript/Reference/Global_Objects/Array"
$ws.Range("B23").Value = $newUrl
$ws.Range("C23").Value = "Métodos para utilizar en arrays javascript"

# Turn B23 into a real hyperlink pointing at the new URL
$ws.Hyperlinks.Add($ws.Range("B23"), $newUrl)

# Re-apply the hyperlink text style (the Hyperlinks.Add call above re-styles
# the cell using a brand new style record); restore it to the same style
# used by the rest of the link column, reusing the existing style record.
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$excel.CutCopyMode = 0
